# Add three new trivia rows (question / 4 answers / correct-answer index)
# to the "hva" sheet, mirroring the two rows already on the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hva")

$ws.Range("A3").Value = "Hvad fuck sker der her"
$ws.Range("B3").Value = "Ingen ved det"
$ws.Range("C3").Value = "Alle vil vide det"
$ws.Range("D3").Value = "Find selv ud af det"
$ws.Range("E3").Value = "Fire jo"
$ws.Range("F3").Value = 2

$ws.Range("A4").Value = "Gustav er fuldstændig uden for pedagogisk rækkeviede?"
$ws.Range("B4").Value = "Yeps"
$ws.Range("C4").Value = "Joooeh"
$ws.Range("D4").Value = "Helt 100"
$ws.Range("E4").Value = "Selvfølgelig"
$ws.Range("F4").Value = 3

$ws.Range("A5").Value = "Hvor meget vejer det tungeste gram?"
$ws.Range("B5").Value = "Mindst 5 ihvertfald"
$ws.Range("C5").Value = 16
$ws.Range("D5").Value = "måske et helt"
$ws.Range("E5").Value = "omtrænt ligeså meget som min højre nos"
$ws.Range("F5").Value = 4

# Re-fit the columns now that column E holds much longer text than before
# (these are the "best fit" widths Excel settled on for the new content).
$ws.Range("A1").ColumnWidth = 55.936
$ws.Range("B1").ColumnWidth = 23.603
$ws.Range("C1").ColumnWidth = 20.926
$ws.Range("D1").ColumnWidth = 21.593
$ws.Range("E1").ColumnWidth = 37.767
$ws.Range("F1").ColumnWidth = 16.593

# Leave the cursor where the author last left it.
$ws.Range("F10").Select()
